$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: given a Word.Paragraph, return its exact original <w:p ...> opening
# tag (attributes like w14:paraId/w:rsid* preserved) by reading WordOpenXML.
# ---------------------------------------------------------------------------
function Get-ParaOpenTag($para) {
    $oxml = $para.Range.WordOpenXML
    if ($oxml -match '(?s)<w:body>(<w:p[^>]*>)') {
        return $matches[1]
    }
    return "<w:p>"
}

# Wrap a body fragment (one or more <w:p>...</w:p> elements) in the minimal
# pkg:package envelope that Range.InsertXML expects.
function Wrap-Body($bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ===========================================================================
# Edit 1: flesh out the empty bullet that followed
# "Sequenz von Events, welche zu einem Deadlock führen" with the
# "Deadlock Pattern der größe k ist:" sub-list.
# ===========================================================================
$rng = $d.Content
$null = $rng.Find.Execute("Sequenz von Events, welche zu einem Deadlock führen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $rng.Paragraphs(1)
$targetPara = $anchorPara.Next()

$openTag = Get-ParaOpenTag $targetPara

$firstPara = $openTag +
    '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Deadlock Pattern der </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>größe</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> k ist:</w:t></w:r>' +
    '</w:p>'

$ilvl1PPr = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'

$p2 = '<w:p>' + $ilvl1PPr + '<w:r><w:t>Sequenz D aus Events</w:t></w:r></w:p>'
$p3 = '<w:p>' + $ilvl1PPr + '<w:r><w:t>Mit k unterschiedlichen Threads</w:t></w:r></w:p>'
$p4 = '<w:p>' + $ilvl1PPr + '<w:r><w:t>Und k unterschiedlichen Locks</w:t></w:r></w:p>'

$p5PPr = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>'
$p5Runs =
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Sodass</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> thread(event) </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>gleich</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> dem</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> T</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">hread </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ist</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>indem</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> das Event </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>läuft</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
$p5 = '<w:p>' + $p5PPr + $p5Runs + '</w:p>'

$p6Runs =
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Op</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>event</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>) ist die O</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">peration des Events, diese muss immer ein </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>acq</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> auf einen Lock sein</w:t></w:r>'
$p6 = '<w:p>' + $ilvl1PPr + $p6Runs + '</w:p>'

$p7 = '<w:p>' + $ilvl1PPr + '</w:p>'

$bodyFragment = $firstPara + $p2 + $p3 + $p4 + $p5 + $p6 + $p7
$targetPara.Range.InsertXML( (Wrap-Body $bodyFragment) )

# ===========================================================================
# Edit 2: add <w:lastRenderedPageBreak/> to the "4.1 Synchronization-
# Preserving Deadlocks" heading run (right before its <w:t>).
# ===========================================================================
$rng = $d.Content
$null = $rng.Find.Execute("4.1 Synchronization-Preserving Deadlocks", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingPara = $rng.Paragraphs(1)

$oxml = $headingPara.Range.WordOpenXML
$null = $oxml -match '(?s)<w:body>(<w:p[ >].*?</w:p>)'
$headingFrag = $matches[1]

$headingFrag = $headingFrag.Replace(
    '<w:t>4.1 Synchronization-Preserving Deadlocks</w:t>',
    '<w:lastRenderedPageBreak/><w:t>4.1 Synchronization-Preserving Deadlocks</w:t>'
)

$headingPara.Range.InsertXML( (Wrap-Body $headingFrag) )

# ===========================================================================
# Edit 3: remove the (now redundant) <w:lastRenderedPageBreak/> from the
# "Das heißt auch, dass andere kritische Sektionen..." run.
# ===========================================================================
$rng = $d.Content
$null = $rng.Find.Execute("Das heißt auch, dass andere kritische Sektionen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$heisstPara = $rng.Paragraphs(1)

$oxml = $heisstPara.Range.WordOpenXML
$null = $oxml -match '(?s)<w:body>(<w:p[ >].*?</w:p>)'
$heisstFrag = $matches[1]

$heisstFrag = $heisstFrag.Replace('<w:lastRenderedPageBreak/>', '')

$heisstPara.Range.InsertXML( (Wrap-Body $heisstFrag) )

Write-Output "All edits applied."
